$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

$ws.Range("C1").Value = "dc.title"
$ws.Range("D1").Value = "dc.type"
$ws.Range("E1").Value = "dc.date.issued"

$ws.Range("B2").Value = "ADD"
$ws.Range("C2").Value = "Test Publication"
$ws.Range("D2").Value = "Article"
$d = Get-Date -Year 2020 -Month 1 -Day 1 -Hour 0 -Minute 0 -Second 0
$ws.Range("E2").Value = $d
$ws.Range("E2").NumberFormat = "yyyy-mm-dd"

$ws.Columns.Item(3).ColumnWidth = 17.333333333333332
$ws.Columns.Item(4).ColumnWidth = 15.0
$ws.Columns.Item(5).ColumnWidth = 14.0

[void]$ws.Activate()
[void]$ws.Range("E3").Select()

Write-Host "done"
